# Commit: gjorde "na" til "NA"  (changed "na" to "NA")
# The lowercase shared string "na" (used as a placeholder for missing
# lactate measurements) is replaced everywhere with uppercase "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$changed = 0
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ceq "na") {
            $cell.Value = "NA"
            $changed = $changed + 1
        }
    }
}

Write-Output "Replaced $changed cell(s) containing 'na' with 'NA'"
